# Rename header columns to use the respective input-file-name suffixes
# (_old -> _FV2410, _new -> _FV2504) instead of the generic _old/_new,
# then wrap the data range in an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into a proper Excel Table (ListObject) with autofilter.
$dataRange = $ws.Range("A1:U81")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)
